$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the explanatory subtitle row entirely (old row 2,
# "(according to the population census data)"); the following empty
# spacer row shifts up to become the new blank row 2.
$ws.Rows.Item(2).Delete()

# Remove the historic 1989/2002 data columns (old columns B and C), keeping the 2014 column
$ws.Range("B1:C1").EntireColumn.Delete()

# Set consistent row heights to match the simplified single-year layout
$ws.Rows.Item(1).RowHeight = 20.1
$ws.Rows.Item(2).RowHeight = 20.1
$ws.Rows.Item(3).RowHeight = 20.1
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
